$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2005420054200542
$ws.Range("C2").Value = 0.5582655826558266
$ws.Range("J2").Value = 0.01084010840108401
$ws.Range("P2").Value = 0.1273712737127371
$ws.Range("S2").Value = 0.1029810298102981
$ws.Range("B3").Value = 0.01415094339622642
$ws.Range("C3").Value = 0.02830188679245283
$ws.Range("J3").Value = 0.04716981132075472
$ws.Range("P3").Value = 0.7216981132075472
$ws.Range("S3").Value = 0.1886792452830189
$ws.Range("J4").Value = 0.02857142857142857
$ws.Range("P4").Value = 0.6285714285714286
$ws.Range("S4").Value = 0.3428571428571429
$ws.Range("B6").Value = 0.05932203389830509
$ws.Range("D6").Value = 0.00423728813559322
$ws.Range("F6").Value = 0.02542372881355932
$ws.Range("J6").Value = 0.3432203389830508
$ws.Range("O6").Value = 0.01694915254237288
$ws.Range("Q6").Value = 0.1016949152542373
$ws.Range("R6").Value = 0.1059322033898305
$ws.Range("S6").Value = 0.3432203389830508
$ws.Range("B7").Value = 0.1063829787234043
$ws.Range("D7").Value = 0.03191489361702127
$ws.Range("F7").Value = 0.05319148936170213
$ws.Range("J7").Value = 0.1382978723404255
$ws.Range("O7").Value = 0.02127659574468085
$ws.Range("Q7").Value = 0.2127659574468085
$ws.Range("R7").Value = 0.04787234042553191
$ws.Range("S7").Value = 0.3882978723404255
$ws.Range("B8").Value = 0.1425061425061425
$ws.Range("D8").Value = 0.01228501228501228
$ws.Range("E8").Value = 0.002457002457002457
$ws.Range("F8").Value = 0.07125307125307126
$ws.Range("J8").Value = 0.1056511056511057
$ws.Range("O8").Value = 0.02457002457002457
$ws.Range("Q8").Value = 0.1572481572481572
$ws.Range("R8").Value = 0.07371007371007371
$ws.Range("S8").Value = 0.4103194103194103
$ws.Range("B9").Value = 0.07926829268292683
$ws.Range("D9").Value = 0.006097560975609756
$ws.Range("F9").Value = 0.06097560975609756
$ws.Range("J9").Value = 0.07317073170731707
$ws.Range("O9").Value = 0.02439024390243903
$ws.Range("Q9").Value = 0.25
$ws.Range("R9").Value = 0.08536585365853659
$ws.Range("S9").Value = 0.4207317073170732
$ws.Range("B10").Value = 0.1352154531946508
$ws.Range("D10").Value = 0.01634472511144131
$ws.Range("F10").Value = 0.06983655274888559
$ws.Range("J10").Value = 0.1367013372956909
$ws.Range("O10").Value = 0.02005943536404161
$ws.Range("Q10").Value = 0.1991084695393759
$ws.Range("R10").Value = 0.0824665676077266
$ws.Range("S10").Value = 0.3402674591381872
$ws.Range("G11").Value = 0.1511627906976744
$ws.Range("J11").Value = 0.1511627906976744
$ws.Range("K11").Value = 0.2412790697674419
$ws.Range("L11").Value = 0.4476744186046512
$ws.Range("S11").Value = 0.008720930232558139
$ws.Range("G12").Value = 0.66875
$ws.Range("J12").Value = 0.25
$ws.Range("K12").Value = 0.00625
$ws.Range("L12").Value = 0.03125
$ws.Range("S12").Value = 0.04375
$ws.Range("G13").Value = 0.6923076923076923
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.0576923076923077
$ws.Range("F15").Value = 0.03139013452914798
$ws.Range("H15").Value = 0.1479820627802691
$ws.Range("I15").Value = 0.06278026905829596
$ws.Range("J15").Value = 0.3677130044843049
$ws.Range("K15").Value = 0.08071748878923767
$ws.Range("M15").Value = 0.0179372197309417
$ws.Range("O15").Value = 0.07623318385650224
$ws.Range("S15").Value = 0.2152466367713005
$ws.Range("F16").Value = 0.02777777777777778
$ws.Range("H16").Value = 0.125
$ws.Range("I16").Value = 0.07870370370370371
$ws.Range("J16").Value = 0.3981481481481481
$ws.Range("K16").Value = 0.1157407407407407
$ws.Range("M16").Value = 0.01851851851851852
$ws.Range("O16").Value = 0.07870370370370371
$ws.Range("S16").Value = 0.1574074074074074
$ws.Range("F17").Value = 0.01847575057736721
$ws.Range("H17").Value = 0.1547344110854504
$ws.Range("I17").Value = 0.08545034642032333
$ws.Range("J17").Value = 0.4110854503464203
$ws.Range("K17").Value = 0.138568129330254
$ws.Range("M17").Value = 0.02540415704387991
$ws.Range("N17").Value = 0.002309468822170901
$ws.Range("O17").Value = 0.05773672055427252
$ws.Range("S17").Value = 0.1062355658198614
$ws.Range("F18").Value = 0.0425531914893617
$ws.Range("H18").Value = 0.1861702127659574
$ws.Range("I18").Value = 0.06914893617021277
$ws.Range("J18").Value = 0.4202127659574468
$ws.Range("K18").Value = 0.1170212765957447
$ws.Range("M18").Value = 0.01063829787234043
$ws.Range("O18").Value = 0.06382978723404255
$ws.Range("S18").Value = 0.09042553191489362
$ws.Range("F19").Value = 0.0271828665568369
$ws.Range("H19").Value = 0.2026359143327842
$ws.Range("I19").Value = 0.06836902800658978
$ws.Range("J19").Value = 0.3887973640856672
$ws.Range("K19").Value = 0.1095551894563427
$ws.Range("M19").Value = 0.02800658978583196
$ws.Range("O19").Value = 0.06589785831960461
$ws.Range("S19").Value = 0.1095551894563427

Write-Output "Applied 108 cell updates"
